$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 50339.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 50339.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 50339.5
$ws.Range("N134").Value = -60479.5
# Row 136
$ws.Range("H136").Value = 54873.332
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 54873.332
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 54873.332
$ws.Range("N136").Value = -65073.332
# Row 139
$ws.Range("H139").Value = 45353.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45353.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45353.332
$ws.Range("N139").Value = -55633.332
# Row 140
$ws.Range("H140").Value = 269800
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 269800
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 269800
$ws.Range("N140").Value = -280160
# Row 141
$ws.Range("H141").Value = 5607.577
$ws.Range("I141").Value = 5659.0415
$ws.Range("J141").Value = 4990
$ws.Range("K141").Value = 16977.1245
$ws.Range("L141").Value = 14970
$ws.Range("M141").Value = -11797.1245
$ws.Range("N141").Value = -25330

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4979.4
$ws.Range("I32").Value = 3879.8909
$ws.Range("J32").Value = 8003.05
$ws.Range("K32").Value = 3879.8909
$ws.Range("L32").Value = 8003.05
$ws.Range("M32").Value = -3592.8909
# Row 61
$ws.Range("H61").Value = 910.95
$ws.Range("I61").Value = 910.95
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 910.95
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -698.95
$ws.Range("N61").ClearContents()
# Row 132
$ws.Range("H132").Value = 2759.95
$ws.Range("I132").Value = 1118.2727
$ws.Range("J132").Value = 4766.4443
$ws.Range("K132").Value = 3354.8181
$ws.Range("L132").Value = 14299.3329
$ws.Range("M132").Value = -824.8181
$ws.Range("N132").Value = -19359.3329
# Row 136
$ws.Range("H136").Value = 910.95
$ws.Range("I136").Value = 910.95
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2732.85
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -182.8500000000004
$ws.Range("N136").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2229.6
$ws.Range("I86").Value = 1944.4
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 1944.4
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -821.4000000000001
$ws.Range("N86").Value = -5046
# Row 89
$ws.Range("H89").Value = 2229.6
$ws.Range("I89").Value = 1944.4
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 9722
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -4106
$ws.Range("N89").Value = -25232
# Row 134
$ws.Range("H134").Value = 3012.8774
$ws.Range("I134").Value = 1134.7715
$ws.Range("J134").Value = 7708.143
$ws.Range("K134").Value = 3404.3145
$ws.Range("L134").Value = 23124.429
$ws.Range("M134").Value = -869.3145000000004
$ws.Range("N134").Value = -28194.429

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 563.7041
$ws.Range("I31").Value = 613.3570999999999
$ws.Range("J31").Value = 543.84283
$ws.Range("K31").Value = 613.3570999999999
$ws.Range("L31").Value = 543.84283
$ws.Range("M31").Value = -318.3570999999999
$ws.Range("N31").Value = -1133.84283
# Row 34
$ws.Range("H34").Value = 563.7041
$ws.Range("I34").Value = 613.3570999999999
$ws.Range("J34").Value = 543.84283
$ws.Range("K34").Value = 613.3570999999999
$ws.Range("L34").Value = 543.84283
$ws.Range("M34").Value = -411.3570999999999
$ws.Range("N34").Value = -947.84283
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 115
$ws.Range("H115").Value = 30825
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 30825
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 30825
$ws.Range("N115").Value = -33175

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1500.4615
$ws.Range("I68").Value = 1045.5
$ws.Range("J68").Value = 1784.8125
$ws.Range("K68").Value = 3136.5
$ws.Range("L68").Value = 5354.4375
$ws.Range("M68").Value = -2325.5
$ws.Range("N68").Value = -6976.4375
# Row 71
$ws.Range("H71").Value = 1500.4615
$ws.Range("I71").Value = 1045.5
$ws.Range("J71").Value = 1784.8125
$ws.Range("K71").Value = 9409.5
$ws.Range("L71").Value = 16063.3125
$ws.Range("M71").Value = -5353.5
$ws.Range("N71").Value = -24175.3125
# Row 76
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 9000
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -9766
# Row 79
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -11652
# Row 94
$ws.Range("H94").Value = 2814.8333
$ws.Range("I94").Value = 945
$ws.Range("J94").Value = 3749.75
$ws.Range("K94").Value = 2835
$ws.Range("L94").Value = 11249.25
$ws.Range("M94").Value = -2159
$ws.Range("N94").Value = -12601.25
# Row 100
$ws.Range("H100").Value = 3964.8
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3964.8
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 11894.4
$ws.Range("N100").Value = -13516.4
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
# Row 106
$ws.Range("H106").Value = 4990
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4990
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 14970
$ws.Range("N106").Value = -16862
# Row 109
$ws.Range("H109").Value = 1870.5
$ws.Range("I109").Value = 744
$ws.Range("J109").Value = 5250
$ws.Range("K109").Value = 2232
$ws.Range("L109").Value = 15750
$ws.Range("M109").Value = -1192
$ws.Range("N109").Value = -17830
# Row 112
$ws.Range("H112").Value = 2027
$ws.Range("I112").Value = 2027
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 6081
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -4973

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3859.1428
$ws.Range("I132").Value = 2127.6155
$ws.Range("J132").Value = 5359.8
$ws.Range("K132").Value = 6382.8465
$ws.Range("L132").Value = 16079.4
$ws.Range("M132").Value = -3852.8465
$ws.Range("N132").Value = -21139.4

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2389
$ws.Range("I100").Value = 2150.5
$ws.Range("J100").Value = 2593.4285
$ws.Range("K100").Value = 2150.5
$ws.Range("L100").Value = 2593.4285
$ws.Range("M100").Value = -1609.5
# Row 132
$ws.Range("H132").Value = 4719.7144
$ws.Range("I132").Value = 3767.5
$ws.Range("J132").Value = 10433
$ws.Range("K132").Value = 11302.5
$ws.Range("L132").Value = 31299
$ws.Range("M132").Value = -8772.5
$ws.Range("N132").Value = -36359
# Row 136
$ws.Range("H136").Value = 3109.851
$ws.Range("I136").Value = 1372.4412
$ws.Range("J136").Value = 7653.846
$ws.Range("K136").Value = 4117.3236
$ws.Range("L136").Value = 22961.538
$ws.Range("M136").Value = -1567.3236
$ws.Range("N136").Value = -28061.538

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2120.3076
$ws.Range("I132").Value = 1433.3572
$ws.Range("J132").Value = 3868.9092
$ws.Range("K132").Value = 4300.071599999999
$ws.Range("L132").Value = 11606.7276
$ws.Range("M132").Value = -1770.071599999999
